# Adds a "NewCustomer" worksheet (with a sample customer row used by the
# "read the excel file" test) after the existing "SearchAddCustomer" sheet,
# and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- add the new worksheet at the end of the tab strip -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "NewCustomer"

# --- row 1: sample customer record ---------------------------------------
$ws.Range("A1").Value = "Adithya"
$ws.Range("B1").Value = "athreya"
$ws.Range("C1").Value = "bengaluru"
$ws.Range("D1").Value = "srinivasa nagar"
$ws.Range("E1").Value = 131
$ws.Range("F1").Value = "adithya.athreya@gmail.com"
$ws.Range("G1").Value = 9742315935
$ws.Range("H1").Value = "This is first customer in POM"

# Email column gets a live mailto: hyperlink with the workbook's normal
# Hyperlink cell style (matches the style already used on the Login sheet).
$ws.Hyperlinks.Add($ws.Range("F1"), "mailto:adithya.athreya@gmail.com") | Out-Null
$ws.Range("F1").Style = "Hyperlink"

# --- best-fit-style column widths for the text/number columns ------------
# (ColumnWidth is expressed in characters; values chosen so the exported
# column width in the saved file lines up with Excel's own "AutoFit"
# result for this content.)
$ws.Columns.Item(3).ColumnWidth = 9.1666666666667
$ws.Columns.Item(4).ColumnWidth = 13.3333333333333
$ws.Columns.Item(6).ColumnWidth = 25.6666666666667
$ws.Columns.Item(7).ColumnWidth = 10.1666666666667
$ws.Columns.Item(8).ColumnWidth = 25.6666666666667

# --- selection / active tab -----------------------------------------------
# Selecting a cell on the new sheet makes it the active sheet (and the
# active tab of the workbook), just like clicking it in the UI.
$ws.Range("H2").Select() | Out-Null
